# Estudio de Roles y Seguridad 3.0 - update IP filter values and view state
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# The "Setting" row values used to restrict by subnet 192.168.3.% are now
# restricted to subnet 192.168.1.% (all J-column "Setting" rows in both
# Paciente and Medico permission blocks).
$ws.Cells.Replace("192.168.3.%", "192.168.1.%", -4142, 1, $false, $false, $false, $false)

# J39 and J58 ("Setting" rows that close each permissions block) pick up the
# same visual style already used by the other IP-restriction cells in their
# column (e.g. J44 / J56) instead of the heavier block-closing style.
$ws.Range("J44").Copy()
$ws.Range("J39").PasteSpecial(-4122)
$ws.Range("J56").Copy()
$ws.Range("J58").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active view: move the selection to I80 (this also drops the
# stale frozen top-left scroll position left over from the previous edit).
$ws.Range("I80").Select()
